$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 160   # was 158
$ws1.Range("F3").Value = 285   # was 283
$ws1.Range("F4").Value = 122   # was 120
$ws1.Range("F5").Value = 1263   # was 1261
$ws1.Range("F6").Value = 17754   # was 17688
$ws1.Range("F7").Value = 341   # was 336
$ws1.Range("F8").Value = 242   # was 235
$ws1.Range("F9").Value = 1061   # was 1060
$ws1.Range("F10").Value = 6709   # was 6697
$ws1.Range("F12").Value = 150   # was 148
$ws1.Range("F14").Value = 103   # was 102
$ws1.Range("F18").Value = 1295   # was 1294
$ws1.Range("F19").Value = 180   # was 175
$ws1.Range("F24").Value = 32   # was 31
$ws1.Range("F25").Value = 258   # was 256
$ws1.Range("F26").Value = 963   # was 957
$ws1.Range("F27").Value = 103   # was 101
$ws1.Range("F28").Value = 5131   # was 5128
$ws1.Range("F29").Value = 530   # was 529
$ws1.Range("F30").Value = 63   # was 61
$ws1.Range("F31").Value = 11893   # was 11868
$ws1.Range("F33").Value = 35   # was 34
$ws1.Range("F34").Value = 194   # was 191
$ws1.Range("F36").Value = 3903   # was 3902

# Sheet "全部类型" (sheet4) - update column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 160   # was 158
$ws4.Range("F3").Value = 285   # was 283
$ws4.Range("F4").Value = 122   # was 120
$ws4.Range("F5").Value = 1263   # was 1261
$ws4.Range("F6").Value = 17754   # was 17688
$ws4.Range("F7").Value = 341   # was 336
$ws4.Range("F8").Value = 242   # was 235
$ws4.Range("F9").Value = 1061   # was 1060
$ws4.Range("F10").Value = 6709   # was 6697
$ws4.Range("F12").Value = 150   # was 148
$ws4.Range("F14").Value = 103   # was 102
$ws4.Range("F18").Value = 1295   # was 1294
$ws4.Range("F19").Value = 180   # was 175
$ws4.Range("F24").Value = 32   # was 31
$ws4.Range("F25").Value = 258   # was 256
$ws4.Range("F26").Value = 963   # was 957
$ws4.Range("F27").Value = 103   # was 101
$ws4.Range("F28").Value = 5131   # was 5128
$ws4.Range("F29").Value = 530   # was 529
$ws4.Range("F32").Value = 63   # was 61
$ws4.Range("F33").Value = 11893   # was 11868
$ws4.Range("F35").Value = 35   # was 34
$ws4.Range("F36").Value = 194   # was 191
$ws4.Range("F38").Value = 3903   # was 3902
